$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $newValue) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $newValue
    $rng.Style = "Normal"
}

# Row 2
Set-TextCell 'D2' '68.022.66'
Set-TextCell 'E2' '  +2.32%  '

# Row 3
Set-TextCell 'D3' '3.617.83'
Set-TextCell 'E3' '  +1.16%  '

# Row 4
Set-TextCell 'D4' '0.998'
Set-TextCell 'E4' '  -0.14%  '

# Row 5
Set-TextCell 'D5' '203.28'
Set-TextCell 'E5' '  +8.44%  '

# Row 6
Set-TextCell 'D6' '563.54'
Set-TextCell 'E6' '  -4.43%  '

# Row 7
Set-TextCell 'D7' '3.617.13'
Set-TextCell 'E7' '  +1.25%  '

# Row 8
Set-TextCell 'D8' '0.616'
Set-TextCell 'E8' '  +0.24%  '

# Row 9
Set-TextCell 'E9' '  +0.19%  '

# Row 10
Set-TextCell 'D10' '0.675'
Set-TextCell 'E10' '  +0.38%  '

# Row 11
Set-TextCell 'D11' '60.36'
Set-TextCell 'E11' '  +12.16%  '

# Row 12
Set-TextCell 'D12' '0.153'
Set-TextCell 'E12' '  +3.86%  '

# Row 13
Set-TextCell 'D13' '0.0000290'
Set-TextCell 'E13' '  +11.42%  '

# Row 14
Set-TextCell 'D14' '10.05'
Set-TextCell 'E14' '  +2.64%  '

# Row 15
Set-TextCell 'D15' '4.203.67'
Set-TextCell 'E15' '  +1.43%  '

# Row 16
Set-TextCell 'D16' '3.606.77'
Set-TextCell 'E16' '  +0.82%  '

# Row 17
Set-TextCell 'E17' '  +0.60%  '

# Row 18
Set-TextCell 'D18' '19.09'
Set-TextCell 'E18' '  +4.09%  '

# Row 19
Set-TextCell 'D19' '67.785.28'
Set-TextCell 'E19' '  +2.03%  '

# Row 20
Set-TextCell 'D20' '12.38'
Set-TextCell 'E20' '  +0.93%  '

# Row 21
Set-TextCell 'E21' '  +2.16%  '

# Row 22
Set-TextCell 'D22' '405.11'
Set-TextCell 'E22' '  +1.83%  '

# Row 23
Set-TextCell 'D23' '12.92'
Set-TextCell 'E23' '  +13.05%  '

# Row 24
Set-TextCell 'D24' '4.17'
Set-TextCell 'E24' '  -4.79%  '

# Row 25
Set-TextCell 'D25' '85.50'
Set-TextCell 'E25' '  -0.36%  '

# Row 26
Set-TextCell 'D26' '2.96'
Set-TextCell 'E26' '  +2.23%  '

# Row 27
Set-TextCell 'D27' '12.59'
Set-TextCell 'E27' '  +0.61%  '

# Row 28
Set-TextCell 'D28' '3.89'
Set-TextCell 'E28' '  +9.13%  '

# Row 29
Set-TextCell 'D29' '6.12'
Set-TextCell 'E29' '  +1.49%  '

# Row 30
Set-TextCell 'D30' '8.39'
Set-TextCell 'E30' '  +17.74%  '

# Row 31
Set-TextCell 'D31' '9.44'
Set-TextCell 'E31' '  +4.93%  '

# Row 32
Set-TextCell 'D32' '31.71'
Set-TextCell 'E32' '  +1.84%  '

# Row 33
Set-TextCell 'D33' '684.29'
Set-TextCell 'E33' '  +10.67%  '

# Row 34
Set-TextCell 'D34' '12.22'
Set-TextCell 'E34' '  +0.63%  '

# Row 35
Set-TextCell 'B35' 'Hedera'
Set-TextCell 'C35' 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextCell 'D35' '0.114'
Set-TextCell 'E35' '  +0.60%  '

# Row 36
Set-TextCell 'B36' 'OKB'
Set-TextCell 'C36' 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextCell 'D36' '63.83'
Set-TextCell 'E36' '  +0.52%  '

# Row 37
Set-TextCell 'D37' '42.63'

# Row 38
Set-TextCell 'D38' '0.431'
Set-TextCell 'E38' '  +9.69%  '

# Row 39
Set-TextCell 'D39' '1.00'

# Row 40
Set-TextCell 'D40' '0.0₃0780'
Set-TextCell 'E40' '  +2.06%  '

# Row 41
Set-TextCell 'D41' '3.25'
Set-TextCell 'E41' '  +15.05%  '

# Row 42
Set-TextCell 'D42' '3.248.68'
Set-TextCell 'E42' '  +7.44%  '

# Row 43
Set-TextCell 'E43' '  +3.66%  '

# Row 44
Set-TextCell 'D44' '2.81'
Set-TextCell 'E44' '  +11.20%  '

# Row 45
Set-TextCell 'D45' '3.06'
Set-TextCell 'E45' '  +30.38%  '

# Row 46
Set-TextCell 'D46' '0.995'
Set-TextCell 'E46' '  -0.34%  '

# Row 47
Set-TextCell 'D47' '0.0419'
Set-TextCell 'E47' '  +2.30%  '

# Row 48
Set-TextCell 'D48' '2.74'
Set-TextCell 'E48' '  +10.25%  '

# Row 49
Set-TextCell 'B49' 'ApeXProtocol'
Set-TextCell 'C49' 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
Set-TextCell 'D49' '3.10'
Set-TextCell 'E49' '  +3.47%  '

# Row 50
Set-TextCell 'B50' 'THORChain'
Set-TextCell 'C50' 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
Set-TextCell 'D50' '8.87'
Set-TextCell 'E50' '  +2.99%  '

# Row 51
Set-TextCell 'B51' 'Stellar'
Set-TextCell 'C51' 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextCell 'D51' '0.131'
Set-TextCell 'E51' '  +0.77%  '
